# Fruta / hortaliza, semanal
# Insert a new weekly record as row 8 (pushing existing rows 8..78 down to 9..79)
# for "Macroferia Regional de Talca" - Arándano (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8; everything below shifts down one row.
$ws.Rows("8:8").Insert()

# Populate the new row 8 with the new data point.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = [DateTime]"2022-11-16"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101001
$ws.Range("J8").Value = "Arándano (blue)"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 3500
$ws.Range("T8").Value = 2
